$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block fix (rows 2 & 3 had swapped/misplaced values) ---
# Row 2: Name: | Curriculum Year: | 2022 (numeric)
$ws.Range("E2").Value = 2022

# Row 3: Student Number: | Department: | Information Technology
$ws.Range("A3").Value = "Student Number:"
$ws.Range("D3").Value = "Department:"
$ws.Range("E3").Value = "Information Technology"

# --- Subject table (rows 6-14) replaced with a new subject list ---
$ws.Range("A6").Value = 25
$ws.Range("B6").Value = "CC 4"
$ws.Range("C6").Value = "Data Structures and Algorithms"
$ws.Range("D6").Value = 3

$ws.Range("A7").Value = 26
$ws.Range("B7").Value = "DIGITAL 1"
$ws.Range("C7").Value = "Digital Logic Design"
$ws.Range("D7").Value = 3

$ws.Range("A8").Value = 28
$ws.Range("B8").Value = "FIL1"
$ws.Range("C8").Value = "Pagsasaling Wika"
$ws.Range("D8").Value = 3

$ws.Range("A9").Value = 29
$ws.Range("B9").Value = "IM 101"
$ws.Range("C9").Value = "Fundamentals of Database Systems"
$ws.Range("D9").Value = 2

$ws.Range("A10").Value = 36
$ws.Range("B10").Value = "NET 101"
$ws.Range("C10").Value = "Networking 1"
$ws.Range("D10").Value = 2

$ws.Range("A11").Value = 37
$ws.Range("B11").Value = "NET 101L"
$ws.Range("C11").Value = "Networking 1 L"
$ws.Range("D11").Value = 1

$ws.Range("A12").Value = 40
$ws.Range("B12").Value = "PATH FIT 3"
$ws.Range("C12").Value = "INDIVIDUAL AND DUAL SPORTS"
$ws.Range("D12").Value = 2

$ws.Range("A13").Value = 42
$ws.Range("B13").Value = "PF 101"
$ws.Range("C13").Value = "Object Oriented Programming"
$ws.Range("D13").Value = 2

$ws.Range("A14").Value = 43
$ws.Range("B14").Value = "PF 101L"
$ws.Range("C14").Value = "Object Oriented Programming"
$ws.Range("D14").Value = 1

# --- Two new subject rows appended (15 & 16), copying formatting from row 14 ---
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(15).PasteSpecial()
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(16).PasteSpecial()

$ws.Range("A15").Value = 45
$ws.Range("B15").Value = "RIZ"
$ws.Range("C15").Value = "Life and Works of Rizal"
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = " "

$ws.Range("A16").Value = 46
$ws.Range("B16").Value = "SP 101"
$ws.Range("C16").Value = "Social Issues and Professional Practices"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = " "
